$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 97
$prev = $row - 1

# Set values first
$ws.Cells.Item($row, 1).Value = 96
$ws.Cells.Item($row, 2).Value = "denmark"
$ws.Cells.Item($row, 3).Value = "1st-division"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45253.79166666666
$ws.Cells.Item($row, 6).Value = "Fredericia"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "B.93"
$ws.Cells.Item($row, 9).Value = 4
$ws.Cells.Item($row, 10).Value = 1.4
$ws.Cells.Item($row, 11).Value = "15/11/2023 19:42"
$ws.Cells.Item($row, 12).Value = 1.53
$ws.Cells.Item($row, 13).Value = "23/11/2023 18:58"
$ws.Cells.Item($row, 14).Value = 4.83
$ws.Cells.Item($row, 15).Value = "15/11/2023 19:42"
$ws.Cells.Item($row, 16).Value = 4.76
$ws.Cells.Item($row, 17).Value = "23/11/2023 18:58"
$ws.Cells.Item($row, 18).Value = 7.16
$ws.Cells.Item($row, 19).Value = "15/11/2023 19:42"
$ws.Cells.Item($row, 20).Value = 5.41
$ws.Cells.Item($row, 21).Value = "23/11/2023 18:58"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/denmark/1st-division/fredericia-boldklubben-1893/pMeVjKZp/"

# Copy cell formatting from the row above for the styled columns (A and E)
# so that the new row reuses the existing cell styles instead of creating new ones.
$ws.Range("A$prev").Copy() | Out-Null
$ws.Range("A$row").PasteSpecial(-4122) | Out-Null

$ws.Range("E$prev").Copy() | Out-Null
$ws.Range("E$row").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
